# Remove the decorative "────…" separator paragraphs and the empty
# spacer paragraphs (<w:p><w:pPr><w:spacing w:before="40"/></w:pPr></w:p>)
# that sit right after every table in the document, while leaving all
# other content (text, tables, images) untouched.

$d = $word.ActiveDocument

# Build the 60-character box-drawing separator text (U+2500) without
# relying on string ctor overloads that may not be available.
$sepChar = [char]9472
$sepText = ""
for ($k = 0; $k -lt 60; $k++) {
    $sepText = $sepText + $sepChar
}

$count = $d.Paragraphs.Count
$ranges = New-Object System.Collections.ArrayList

for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text
    $len = $text.Length

    if ($len -ge 1) {
        # Strip the trailing paragraph mark before comparing.
        $core = $text.Substring(0, $len - 1)
    } else {
        $core = $text
    }

    $isSeparator = ($core -eq $sepText)
    $isEmptySpacer = ($core -eq "") -and ($para.Range.ParagraphFormat.SpaceBefore -eq 2)

    if ($isSeparator -or $isEmptySpacer) {
        [void]$ranges.Add("$($para.Range.Start)-$($para.Range.End)")
    }
}

# Delete from the last match to the first so earlier offsets stay valid.
for ($k = $ranges.Count - 1; $k -ge 0; $k--) {
    $parts = $ranges[$k].Split("-")
    $s = [int]$parts[0]
    $e = [int]$parts[1]
    $rng = $d.Range($s, $e)
    $rng.Delete()
}

Write-Output "Removed $($ranges.Count) paragraphs (separators + empty spacers)."
